$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> [new nombre_aides (C), new montant_total (E)]
$updates = @{
    8   = @(181370, 653248410)
    10  = @(278210, 1752513740)
    19  = @(108926, 344804726)
    97  = @(98510, 307118902)
    152 = @(126054, 716092217)
    164 = @(50586, 168946548)
    168 = @(285122, 1213790875)
    169 = @(562673, 1286266968)
    170 = @(367577, 2848223671)
    171 = @(115230, 449025067)
    174 = @(357382, 1020222665)
    175 = @(125700, 815905290)
    179 = @(235806, 813703438)
    180 = @(141532, 341248758)
    279 = @(28968, 57089102)
    293 = @(61673, 194903805)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

$wb.Save()
